$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOMPlanTest")

# New import file and parameters: add a new column L value mirroring E3 (46.875)
$ws.Range("L3").Value = 46.875

# Data valid check for each table: quantity in row 7 (SA-Test3) bumped from 3 to 6
$ws.Range("F7").Value = 6

# Remove the explicit style on E3 so it reverts to the default (general) style
$ws.Range("E3").Style = "Normal"

# Update the active selection to reflect the cell the user was last working with
$ws.Activate()
$ws.Range("F9").Select()
